$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite rows 2-5 with the 2010-2013 data (previously rows 8-11),
# leaving the header row (row 1) untouched.
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 331.441152
$ws.Range("C2").Value = 533.264832
$ws.Range("D2").Value = 201.82368

$ws.Range("A3").Value = "2011年"
$ws.Range("B3").Value = 425.675023
$ws.Range("C3").Value = 698.426139
$ws.Range("D3").Value = 272.751116

$ws.Range("A4").Value = "2012年"
$ws.Range("B4").Value = 512.981851
$ws.Range("C4").Value = 821.121678
$ws.Range("D4").Value = 308.139827

$ws.Range("A5").Value = "2013年"
$ws.Range("B5").Value = 608.77119814
$ws.Range("C5").Value = 922.75239655
$ws.Range("D5").Value = 313.98119841

# Remove the old rows 6-11 (years 2008-2009 leftovers plus the
# duplicated 2010-2013 rows that used to live there) entirely so the
# used range shrinks back down to A1:D5.
$ws.Range("A6:D11").EntireRow.Delete()
